$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 173 (API-Catalogue "undelete" section is alphabetically
# sorted; the new entry "setBusinessDocumentNumbering" belongs right after
# "setBusinessDocument" (row 172) and before "setBusinessDocumentType"
# (old row 173), so everything from old row 173 downward shifts down by one).
$ws.Rows(173).Insert()

# The freshly inserted row has no formatting of its own - copy the
# (now shifted-down) row 174's formats onto it so it matches the rest of
# the table (borders/fill/font as used by every other data row).
$ws.Range("B174:C174").Copy()
$ws.Range("B173:C173").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new API entry.
$ws.Range("B173").Value = "transaction.undelete.master.setBusinessDocumentNumbering"
$ws.Range("C173").Value = "Membatalkan Penghapusan Data Penomoran Dokumen Bisnis"

# Update the view state left by the editor: scrolled/frozen pane anchored
# further down the (now longer) list, with C39 -> B176 as the active cell.
$ws.Range("B176").Select()
